# ==========================================================================
# Weekly CompStat (28th Precinct) refresh: new volume/week header + updated
# crime-complaint figures for rows 15-31 (Rape .. Hate Crimes).
# ==========================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and reporting week -----------------------
$ws.Range("A8").Value = "Volume 32   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  3/31/2025  Through  4/6/2025"

# --- Crime-complaint figures (rows 15-31) ----------------------------------
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -42.857142857142

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -87.5
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -68.75
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = -35.294117647058
$ws.Range("L16").Value = -23.255813953488
$ws.Range("M16").Value = -37.735849056603
$ws.Range("N16").Value = -84.792626728110

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -21.739130434782
$ws.Range("I17").Value = 52
$ws.Range("J17").Value = 66
$ws.Range("K17").Value = -21.212121212121
$ws.Range("L17").Value = -25.714285714285
$ws.Range("M17").Value = -3.703703703703
$ws.Range("N17").Value = -65.562913907284

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 8
$ws.Range("I18").Value = 22
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -24.137931034482
$ws.Range("L18").Value = -37.142857142857
$ws.Range("M18").Value = -18.518518518518
$ws.Range("N18").Value = -90.557939914163

# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -22.222222222222
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -40.909090909090
$ws.Range("I19").Value = 74
$ws.Range("J19").Value = 132
$ws.Range("K19").Value = -43.939393939393
$ws.Range("L19").Value = -25.252525252525
$ws.Range("M19").Value = -2.631578947368
$ws.Range("N19").Value = -29.523809523809

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "'0"
$ws.Range("E20").Value = "'***.*"
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 9
$ws.Range("K20").Value = -35.714285714285
$ws.Range("L20").Value = -40
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -80.434782608695

# Row 21
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = -35.416666666666
$ws.Range("I21").Value = 195
$ws.Range("J21").Value = 297
$ws.Range("K21").Value = -34.343434343434
$ws.Range("L21").Value = -26.691729323308
$ws.Range("M21").Value = -9.302325581395
$ws.Range("N21").Value = -74.708171206225

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = -42.857142857142
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = -33.333333333333

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = -23.076923076923
$ws.Range("L23").Value = -23.076923076923
$ws.Range("M23").Value = 122.222222222222

# Row 24
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -28.947368421052
$ws.Range("F24").Value = 125
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = 16.822429906542
$ws.Range("I24").Value = 425
$ws.Range("J24").Value = 288
$ws.Range("K24").Value = 47.569444444444
$ws.Range("L24").Value = 39.344262295082
$ws.Range("M24").Value = 88.053097345132

# Row 25
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -39.285714285714
$ws.Range("F25").Value = 69
$ws.Range("G25").Value = 77
$ws.Range("H25").Value = -10.389610389610
$ws.Range("I25").Value = 284
$ws.Range("J25").Value = 173
$ws.Range("K25").Value = 64.161849710982
$ws.Range("L25").Value = 45.641025641025

# Row 26
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 600
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 115
$ws.Range("J26").Value = 96
$ws.Range("K26").Value = 19.791666666666
$ws.Range("L26").Value = 16.161616161616
$ws.Range("M26").Value = 3.603603603603

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = -20
$ws.Range("L27").Value = 33.333333333333

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 13
$ws.Range("K28").Value = -30.769230769230
$ws.Range("L28").Value = -18.181818181818

# Row 29
$ws.Range("M29").Value = -80
$ws.Range("N29").Value = -96.296296296296

# Row 30
$ws.Range("M30").Value = -80
$ws.Range("N30").Value = -95.652173913043

# Row 31
$ws.Range("L31").Value = -75

# --- Fix up number formats / cell styles where a cell switched between ----
# --- a literal dash ("0"/"***.*") placeholder and a real number --------
$excel.CutCopyMode = $false

$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
